# Added graticuleps and plotps
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New Notes (column C) content for a few rows ---
# Order matters for shared-string index assignment (mirrors authoring order).

# Row 21 (gridded_flux): replace "n/a" with a new note
$ws.Range("C21").Value = "need its_live and measures velocity data first"

# Row 42 (inset_unproj): add a new note
$ws.Range("C42").Value = "think it makes more sense to have mapzoom first"

# Row 40 (coord): add a new note
$ws.Range("C40").Value = "interactive (mouse clicks)"

# Row 33 (graticuleps): add a new note
$ws.Range("C33").Value = "looks right??"

# --- Column C width grew to fit the new, longer notes ---
$ws.Columns.Item(3).ColumnWidth = 42.36328125

# --- Selection / scroll position left where the author was working ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
